$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("M1").Value = "Linked"
$ws.Range("M1").Select() | Out-Null
